# "novo modelo de composição" — update header info, replace line-item 1 with a
# new product, and reset line-items 2..23 (rows 12..33) back to the blank
# template state used by the remaining rows further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (rows 3-5) --------------------------------------------
$ws.Range("C3").Value2 = "Luiz Henrique"
$ws.Range("C4").Value2 = "Larissa Sousa"

# C5 currently holds "1020" stored as TEXT even though the cell's number
# format is a plain numeric format. Assigning Value2 with a numeric-looking
# string gets auto-coerced to a real number by Excel, which would change
# both the stored type and (indirectly) invite a style change. Instead,
# write it as a text FORMULA ("="1026"") and then paste-special just the
# resulting value back onto itself: this keeps the result as literal text
# without touching the cell's number format / style id.
$c5 = $ws.Range("C5")
$c5.Formula = '="1026"'
$c5.Copy()
$c5.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Line item 1 (row 11): swap in the new product ----------------------
$ws.Range("B11").Value2 = 1
$ws.Range("D11").Value2 = "010 DIFUSOR VSD 35 4 1200 0 F AN0 TROX"
$ws.Range("G11").Value2 = "VSD35-4F01200x0x00000AN0"
$ws.Range("I11").Value2 = 247.64

# --- Line items 2..23 (rows 12-33): clear back to blank template rows ---
$ws.Range("B12:D33").ClearContents()
$ws.Range("G12:K33").ClearContents()
$ws.Range("E12:E33").Value2 = "-"
